# Wording changes from PM&C.
# Applies to dashboard_loader/indigenous_employment_uploader/indigenous_employment.xlsx
#
# 1. On the "Description" sheet: insert a new note paragraph (row 9) above the
#    existing "Notes" block, add a new "References" row at the bottom, and
#    update the view's active selection.
# 2. Shared strings / cell text changes ride along automatically from the
#    Range.Value assignments below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")
$ws.Activate()

# --- Insert the new explanatory paragraph as a new row 9 (pushes the old
#     rows 9-11 "Notes"/"Source" block down to 10-12). ---
$ws.Rows.Item(9).Insert()

# New paragraph text, styled like the other note paragraphs above it
# (wrapped text, same body font) but with pure black font colour.
$ws.Range("B9").Value = "State and Territory employment outcomes are influenced by economic circumstances both within their jurisdictions and across the nation. It should be noted that the Commonwealth has primary responsibility for national economic management and delivers labour market assistance to jobseekers."

$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").WrapText = $true
$ws.Range("B9").Font.Color = 0
$ws.Range("B9").Font.Size = 12

$ws.Rows.Item(9).RowHeight = 39.55

# Leave A9 empty (style carried over from the row above via the insert).

# --- Append the new "References" row after the existing content. ---
$ws.Range("A13").Value = "References"
$ws.Range("B13").Value = "COAG Reform Council, Indigenous Reform 2012-13: Five years of performance, p. 68."
$ws.Rows.Item(13).RowHeight = 12.8

# --- Update the sheet view: scrolled down a bit, with a new active cell. ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("J20").Select()
